# Updated cryptos list on Mon Mar 27 18:56:19 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) columns of the
# crypto-ranking sheet with freshly scraped values. Column D holds numeric-
# looking strings (e.g. "27.119.51", "1.001") that must stay plain TEXT, just
# like they were before the edit (coinranking's "price" text is dot-grouped
# and would be silently reinterpreted as a number/date by Excel's normal
# auto-detection). Prefixing the literal with a leading apostrophe forces
# Excel to store it as text without touching the cell's number format, and
# ClearFormats() afterwards drops the transient "quote prefix" flag so the
# cell's style stays exactly as it was (no stray Text number-format style
# left behind). Column E ("  -2.69%  " etc.) already contains padding
# whitespace so Excel stores it as text natively - no special handling
# needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "27.119.51"; E = "  -2.69%  " },
    @{ Row = 3; D = "1.717.76"; E = "  -2.93%  " },
    @{ Row = 4; D = "1.001"; E = "  -0.07%  " },
    @{ Row = 5; D = "307.60"; E = "  -6.27%  " },
    @{ Row = 6; D = $null; E = "  -0.10%  " },
    @{ Row = 7; D = "0.4710"; E = "  +5.37%  " },
    @{ Row = 8; D = "0.3427"; E = "  -3.77%  " },
    @{ Row = 9; D = "42.18"; E = "  +0.39%  " },
    @{ Row = 10; D = "0.07271"; E = "  -2.35%  " },
    @{ Row = 11; D = "1.042"; E = "  -5.00%  " },
    @{ Row = 12; D = "0.9999"; E = "  -0.21%  " },
    @{ Row = 13; D = "19.91"; E = "  -5.00%  " },
    @{ Row = 14; D = "5.848"; E = "  -2.94%  " },
    @{ Row = 15; D = "1.717.41"; E = "  -3.06%  " },
    @{ Row = 16; D = "6.887"; E = "  -4.73%  " },
    @{ Row = 17; D = "88.85"; E = "  -4.90%  " },
    @{ Row = 18; D = "0.00001037"; E = "  -2.09%  " },
    @{ Row = 19; D = "0.06354"; E = "  -1.19%  " },
    @{ Row = 20; D = "0.9999"; E = "  -0.11%  " },
    @{ Row = 21; D = "16.49"; E = "  -3.63%  " },
    @{ Row = 22; D = "5.619"; E = "  -2.80%  " },
    @{ Row = 23; D = "27.184.06"; E = "  -2.62%  " },
    @{ Row = 24; D = "10.84"; E = "  -4.02%  " },
    @{ Row = 25; D = "2.145"; E = "  +1.19%  " },
    @{ Row = 26; D = "157.35"; E = "  -3.66%  " },
    @{ Row = 27; D = "19.46"; E = "  -4.50%  " },
    @{ Row = 28; D = "1.911.83"; E = "  -3.21%  " },
    @{ Row = 29; D = "2.098"; E = "  -3.12%  " },
    @{ Row = 30; D = "119.77"; E = "  -4.31%  " },
    @{ Row = 31; D = "1.021"; E = "  -7.53%  " },
    @{ Row = 32; D = "0.09173"; E = "  -0.06%  " },
    @{ Row = 33; D = "3.586"; E = "  -1.88%  " },
    @{ Row = 34; D = "5.319"; E = "  -5.01%  " },
    @{ Row = 35; D = "0.02196"; E = "  -4.23%  " },
    @{ Row = 36; D = "0.05816"; E = "  -5.07%  " },
    @{ Row = 37; D = "10.96"; E = "  -7.74%  " },
    @{ Row = 38; D = "0.1993"; E = "  -5.14%  " },
    @{ Row = 39; D = "4.731"; E = "  -4.67%  " },
    @{ Row = 40; D = "1.393"; E = "  -0.14%  " },
    @{ Row = 41; D = $null; E = "  -6.76%  " },
    @{ Row = 42; D = $null; E = "  -5.74%  " },
    @{ Row = 43; D = "7.435"; E = "  -6.01%  " },
    @{ Row = 44; D = $null; E = "  -5.85%  " },
    @{ Row = 45; D = "0.5638"; E = "  -4.33%  " },
    @{ Row = 46; D = "3.554"; E = "  -5.05%  " },
    @{ Row = 47; D = "117.79"; E = "  -3.73%  " },
    @{ Row = 48; D = $null; E = "  -5.92%  " },
    @{ Row = 49; D = $null; E = "  -3.75%  " },
    @{ Row = 50; D = "1.085"; E = "  -4.52%  " },
    @{ Row = 51; D = "0.9997"; E = "  -0.09%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        $ws.Cells.Item($row, 4).Value = "'" + $u.D
        $ws.Cells.Item($row, 4).ClearFormats()
    }
    $ws.Cells.Item($row, 5).Value = $u.E
}
